# Update MSME Country Indicators - Korea, Rep. Summary figures.
# These cells hold their numbers as text (General format, shared strings in
# the OOXML), so a plain numeric assignment would silently convert them to
# real numbers. Prefixing with a single quote forces Excel to keep them as
# text, matching the source data's storage as strings like "66.6" / "72".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) -- row 11
$ws.Range("B11").Value = "'66.59"
$ws.Range("C11").Value = "'5.39"
$ws.Range("D11").Value = "'71.98"

# Employment (% of total) -- row 12
$ws.Range("B12").Value = "'41.14"
$ws.Range("C12").Value = "'44.92"
$ws.Range("D12").Value = "'86.06"

# Enterprises (% of total) -- row 14
$ws.Range("B14").Value = "'92.43"
$ws.Range("C14").Value = "'7.48"
$ws.Range("D14").Value = "'99.91"
